$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 8, shifting existing rows 8-10 down to 9-11
$ws.Rows.Item(8).Insert()

# Copy the date style (s="2") used by column D from the row above into the new row's D cell
$ws.Cells.Item(7, 4).Copy()
$ws.Cells.Item(8, 4).PasteSpecial(-4122)  # xlPasteFormats

# Fill the new row 8 with the new weekly record
$ws.Cells.Item(8, 1).Value = 12
$ws.Cells.Item(8, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(8, 3).Value = "Metropolitana"
$ws.Cells.Item(8, 4).Value = 44449
$ws.Cells.Item(8, 5).Value = 13
$ws.Cells.Item(8, 6).Value = 100112026
$ws.Cells.Item(8, 7).Value = "Haba"
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 30
$ws.Cells.Item(8, 11).Value = 16000
$ws.Cells.Item(8, 12).Value = 16000
$ws.Cells.Item(8, 13).Value = 16000
$ws.Cells.Item(8, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(8, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(8, 16).Value = 640
$ws.Cells.Item(8, 17).Value = 25
$ws.Cells.Item(8, 18).Value = "Hortaliza"
